# Add a new worksheet "TestQueries" at the end of the workbook.
$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tq = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$tq.Name = "TestQueries"

# ---------------------------------------------------------------------------
# Sheet: Student Data
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Student Data")
$ws.Range("G2").Value = "checked addStudent and updateStudentToken at db level"
$ws.Range("G3").Value = "checked addStudent at php level"
$ws.Range("B4:E5").Select()

# ---------------------------------------------------------------------------
# Sheet: Professor Data
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Professor Data")
$ws.Range("F2").Value = "syzygy"

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

$ws.Range("G2").Formula = "=""CALL addProfessor('"" & B2 & ""', '"" & C2 & ""', '"" & D2 & ""', '"" & E2 & ""', '"" & F2 & ""');"""
$ws.Range("G3").Formula = "=""CALL addProfessor('"" & B3 & ""', '"" & C3 & ""', '"" & D3 & ""', '"" & E3 & ""', '"" & F3 & ""');"""
$ws.Range("G4").Formula = "=""CALL addProfessor('"" & B4 & ""', '"" & C4 & ""', '"" & D4 & ""', '"" & E4 & ""', '"" & F4 & ""');"""
$ws.Range("G5").Formula = "=""CALL addProfessor('"" & B5 & ""', '"" & C5 & ""', '"" & D5 & ""', '"" & E5 & ""', '"" & F5 & ""');"""

$ws.Range("E8").Value = "checked addProfessor and updateProfessorToken at database level"
$ws.Range("E9").Value = "checked professor_first_login, getProfessorByLoginid and updateProfessorToken in php"

$ws.Range("B3:F3").Select()

# ---------------------------------------------------------------------------
# Sheet: Course Data
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Course Data")
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

$ws.Range("E2").Formula = "=""CALL addCourse('"" & B2 & ""', '"" & C2 & ""');"""
$ws.Range("E3").Formula = "=""CALL addCourse('"" & B3 & ""', '"" & C3 & ""');"""
$ws.Range("E4").Formula = "=""CALL addCourse('"" & B4 & ""', '"" & C4 & ""');"""
$ws.Range("E5").Formula = "=""CALL addCourse('"" & B5 & ""', '"" & C5 & ""');"""

$ws.Range("J2").Value = "addCourse works at DB level"
$ws.Range("J3").Value = "There is no php code to test for this one"

$ws.Range("A6").Select()

# ---------------------------------------------------------------------------
# Sheet: Section Data
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Section Data")

$bcValues = @(
    @(1,1),
    @(1,1),
    @(1,2),
    @(2,1),
    @(2,1),
    @(2,2),
    @(3,2),
    @(3,1),
    @(3,1),
    @(4,2),
    @(4,2),
    @(4,3)
)
for ($i = 0; $i -lt $bcValues.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $i + 1
    $ws.Range("B$r").Value = $bcValues[$i][0]
    $ws.Range("C$r").Value = $bcValues[$i][1]
    $ws.Range("I$r").Formula = "=""CALL addSection("" & B$r & "", "" & C$r & "", "" & D$r & "");"""
}

$ws.Range("A2:A13").Select()

# ---------------------------------------------------------------------------
# Sheet: TestQueries (new sheet)
# ---------------------------------------------------------------------------
$ws = $tq

$ws.Range("A1").Value = "function name"
$ws.Range("B1").Value = "www.attend-in.com"
$ws.Range("C1").Value = "attend-in.php"

# Block 1: student_first_login
$ws.Range("A2").Value = "student_first_login"
$ws.Range("B2").Value = "loginid"
$ws.Range("C2").Value = "last_name"
$ws.Range("D2").Value = "first_name"
$ws.Range("E2").Value = "token"
$ws.Range("I2").Value = "Query"
$ws.Range("J2").Value = "Expected outcome"

$ws.Range("B3").Value = "ngeyer1"
$ws.Range("C3").Value = "Geyer"
$ws.Range("D3").Value = "Nathan"
$ws.Range("E3").Value = "zyzzyvas"
$ws.Range("I3").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$2 & "&" & $B$2 & "=" & $B3 & "&" & $C$2 & "=" & $C3 & "&" & $D$2 & "=" & $D3 & "&" & $E$2 & "=" & $E3'
$ws.Range("J3").Value = "Add ngeyer1 w/token"

$ws.Range("B4").Value = "cchristine1"
$ws.Range("C4").Value = "Sanders"
$ws.Range("D4").Value = "China"
$ws.Range("E4").Value = "zizzling"
$ws.Range("I4").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$2 & "&" & $B$2 & "=" & $B4 & "&" & $C$2 & "=" & $C4 & "&" & $D$2 & "=" & $D4 & "&" & $E$2 & "=" & $E4'
$ws.Range("J4").Value = "Added cchristine then added token."

$ws.Range("B5").Value = "mkandagadda1"
$ws.Range("C5").Value = "Kandagadda"
$ws.Range("D5").Value = "Mounika"
$ws.Range("E5").Value = "jazzlike"
$ws.Range("I5").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$2 & "&" & $B$2 & "=" & $B5 & "&" & $C$2 & "=" & $C5 & "&" & $D$2 & "=" & $D5 & "&" & $E$2 & "=" & $E5'

$ws.Range("B6").Value = "sfarah1"
$ws.Range("C6").Value = "Farah"
$ws.Range("D6").Value = "Sharmarke"
$ws.Range("E6").Value = "quizzing"
$ws.Range("I6").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$2 & "&" & $B$2 & "=" & $B6 & "&" & $C$2 & "=" & $C6 & "&" & $D$2 & "=" & $D6 & "&" & $E$2 & "=" & $E6'

# Block 2: student_login
$ws.Range("A8").Value = "student_login"
$ws.Range("B8").Value = "hashtime"
$ws.Range("C8").Value = "md5_hash"
$ws.Range("D8").Value = "student_tid"
$ws.Range("E8").Value = "class_tid"
$ws.Range("F8").Value = "ip_address"
$ws.Range("G8").Value = "latitude"
$ws.Range("H8").Value = "longitude"

$ws.Range("I9").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$8 & "&" & $B$8 & "=" & $B9 & "&" & $C$8 & "=" & $C9 & "&" & $D$8 & "=" & $D9 & "&" & $E$8 & "=" & $E9 & "&" & $F$8 & "=" & $F9 & "&" & $G$8 & "=" & $G9 & "&" & $H$8 & "=" & $H9'

# Block 3: student_list
$ws.Range("A14").Value = "student_list"
$ws.Range("B14").Value = "hashtime"
$ws.Range("C14").Value = "md5_hash"
$ws.Range("D14").Value = "student_tid"

$ws.Range("I15").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$14 & "&" & $B$14 & "=" & $B15 & "&" & $C$14 & "=" & $C15 & "&" & $D$14 & "=" & $D15'

# Block 4: professor_first_login
$ws.Range("A20").Value = "professor_first_login"
$ws.Range("B20").Value = "loginid"
$ws.Range("C20").Value = "title"
$ws.Range("D20").Value = "last_name"
$ws.Range("E20").Value = "first_name"
$ws.Range("F20").Value = "token"

$ws.Range("B21").Value = "cfrederick"
$ws.Range("C21").Value = "Professor"
$ws.Range("D21").Value = "Frederick"
$ws.Range("E21").Value = "Chad"
$ws.Range("F21").Value = "bezazzes"
$ws.Range("I21").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$20 & "&" & $B$20 & "=" & $B21 & "&" & $C$20 & "=" & $C21 & "&" & $D$20 & "=" & $D21 & "&" & $E$20 & "=" & $E21 & "&" & $F$20 & "=" & $F21'

# Block 5: professor_class_list
$ws.Range("A26").Value = "professor_class_list"
$ws.Range("B26").Value = "hashtime"
$ws.Range("C26").Value = "md5_hash"
$ws.Range("D26").Value = "professor_tid"

$ws.Range("I27").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$26 & "&" & $B$26 & "=" & $B27 & "&" & $C$26 & "=" & $C27 & "&" & $D$26 & "=" & $D27'

# Block 6: set_class
$ws.Range("A32").Value = "set_class"
$ws.Range("B32").Value = "hashtime"
$ws.Range("C32").Value = "md5_hash"
$ws.Range("D32").Value = "professor_tid"
$ws.Range("E32").Value = "class_tid"

$ws.Range("I33").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$32 & "&" & $B$32 & "=" & $B33 & "&" & $C$32 & "=" & $C33 & "&" & $D$32 & "=" & $D33 & "&" & $E$32 & "=" & $E33'

# Block 7: attendance_by_date
$ws.Range("A38").Value = "attendance_by_date"
$ws.Range("B38").Value = "hashtime"
$ws.Range("C38").Value = "md5_hash"
$ws.Range("D38").Value = "professor_tid"
$ws.Range("E38").Value = "class_tid"
$ws.Range("F38").Value = "class_date"

$ws.Range("I39").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$38 & "&" & $B$38 & "=" & $B39 & "&" & $C$38 & "=" & $C39 & "&" & $D$38 & "=" & $D39 & "&" & $E$38 & "=" & $E39 & "&" & $F$38 & "=" & $F39'

# Block 8: attendance_by_student
$ws.Range("A44").Value = "attendance_by_student"
$ws.Range("B44").Value = "hashtime"
$ws.Range("C44").Value = "md5_hash"
$ws.Range("D44").Value = "professor_tid"
$ws.Range("E44").Value = "class_tid"
$ws.Range("F44").Value = "student_tid"

$ws.Range("I45").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$44 & "&" & $B$44 & "=" & $B45 & "&" & $C$44 & "=" & $C45 & "&" & $D$44 & "=" & $D45 & "&" & $E$44 & "=" & $E45 & "&" & $F$44 & "=" & $F45'

# Block 9: set_attendance
$ws.Range("A50").Value = "set_attendance"
$ws.Range("B50").Value = "hashtime"
$ws.Range("C50").Value = "md5_hash"
$ws.Range("D50").Value = "professor_tid"
$ws.Range("E50").Value = "class_tid"
$ws.Range("F50").Value = "student_tid"
$ws.Range("G50").Value = "class_date"
$ws.Range("H50").Value = "attendance_code"

$ws.Range("I51").Formula = '=$B$1 & "/" & $C$1 & "?function=" & $A$50 & "&" & $B$50 & "=" & $B51 & "&" & $C$50 & "=" & $C51 & "&" & $D$50 & "=" & $D51 & "&" & $E$50 & "=" & $E51 & "&" & $F$50 & "=" & $F51 & "&" & $G$50 & "=" & $G51 & "&" & $H$50 & "=" & $H51'

$ws.Range("H27").Select()

Write-Host "Done with part 1"
